# Finished input file scripts for ref scenario
#
# 1. Sheet "All_Seasons_normal" (sheet2.xml): add a new column X = "First Year"
#    with per-row values, and clear the leftover sort state on A21:G25.
# 2. Sheet "CombinedHydroSolar" (sheet1.xml): change the selection to the
#    whole column O (O1:O1048576).
# The workbook's active sheet stays "All_Seasons_normal", selection F9.

$wb = $excel.ActiveWorkbook

# --- Sheet "CombinedHydroSolar": select whole column O -------------------
$ws1 = $wb.Worksheets.Item("CombinedHydroSolar")
$ws1.Activate()
$ws1.Columns("O").Select()

# --- Sheet "All_Seasons_normal": new "First Year" column + cleanup -------
$ws2 = $wb.Worksheets.Item("All_Seasons_normal")
$ws2.Activate()

$ws2.Range("X1").Value = "First Year"

$firstYear = @{
    2  = 2010
    3  = 2010
    4  = 2010
    5  = 2023
    6  = 2030
    7  = 2026
    8  = 2028
    9  = 2025
    11 = 2009
    12 = 1966
    13 = 1962
    14 = 2003
    15 = 2030
    16 = 2025
    17 = 2024
    18 = 2028
    19 = 2030
    20 = 2030
    21 = 2010
    22 = 2010
    23 = 2018
    24 = 2010
    25 = 2000
}

foreach ($row in $firstYear.Keys) {
    $ws2.Cells.Item($row, 24).Value = $firstYear[$row]
}

# Remove the leftover sortState left over from a previous Data > Sort.
$ws2.Sort.SortFields.Clear()

# Restore the sheet's own selection.
$ws2.Range("F9").Select()
